# Update the "想去人数" (F column) counts across the four sheets to the
# newly generated values.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 21447
$ws1.Range("F3").Value  = 3405
$ws1.Range("F4").Value  = 863
$ws1.Range("F6").Value  = 563
$ws1.Range("F7").Value  = 820
$ws1.Range("F8").Value  = 309
$ws1.Range("F11").Value = 149
$ws1.Range("F12").Value = 598
$ws1.Range("F14").Value = 383
$ws1.Range("F15").Value = 46
$ws1.Range("F17").Value = 256
$ws1.Range("F18").Value = 48
$ws1.Range("F20").Value = 93
$ws1.Range("F21").Value = 169

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value  = 153

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 6188
$ws3.Range("F3").Value = 737
$ws3.Range("F4").Value = 737
$ws3.Range("F5").Value = 1757
$ws3.Range("F6").Value = 102

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 6188
$ws4.Range("F3").Value  = 737
$ws4.Range("F4").Value  = 737
$ws4.Range("F5").Value  = 1757
$ws4.Range("F6").Value  = 21447
$ws4.Range("F7").Value  = 3405
$ws4.Range("F8").Value  = 863
$ws4.Range("F9").Value  = 153
$ws4.Range("F10").Value = 102
$ws4.Range("F12").Value = 563
$ws4.Range("F13").Value = 820
$ws4.Range("F14").Value = 309
$ws4.Range("F20").Value = 149
$ws4.Range("F23").Value = 598
$ws4.Range("F27").Value = 383
$ws4.Range("F29").Value = 46
$ws4.Range("F32").Value = 257
$ws4.Range("F33").Value = 48
$ws4.Range("F37").Value = 93
$ws4.Range("F43").Value = 169
